# Generate Report for Handback
#
# This script reflects a "handback" event for the a.md file: both the
# zh-cn and de-de localization targets now have a handback xliff file and
# a handback datetime, so their Status moves from "Ready for handoff" to
# "Handed back: in sync with en-US" and their "Latest Target File" /
# "Latest Handback File" / "Latest Handback DateTime" columns get filled
# in (b.md is left untouched / still pending).

$wb = $excel.ActiveWorkbook

$newStatus   = "Handed back: in sync with en-US"
$aMdUrl      = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/18e8ee69f6769db6f22f2b653512af00b5313b56/e2e/a.md"
$zhHandback  = "a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.zh-cn.xlf"
$deHandback  = "a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.de-de.xlf"
$zhDateTime  = "2016-12-15 03:59:35"
$deDateTime  = "2016-12-15 03:59:54"

# ---------------------------------------------------------------------
# Overview sheet: update the zh-cn / de-de status summary columns (E, F)
# and widen those columns to fit the new, longer status text.
# ---------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")

$wsOverview.Range("E2").Value = $newStatus
$wsOverview.Range("F2").Value = $newStatus
$wsOverview.Range("E3").Value = $newStatus
$wsOverview.Range("F3").Value = $newStatus

$wsOverview.Columns.Item(5).ColumnWidth = 29.14
$wsOverview.Columns.Item(6).ColumnWidth = 29.14

# ---------------------------------------------------------------------
# zh-cn sheet: record the handback for a.md
# ---------------------------------------------------------------------
$wsZh = $wb.Worksheets.Item("zh-cn")

$wsZh.Range("C2").Value = $newStatus
$wsZh.Range("C3").Value = $newStatus

# "Latest Target File" (J) now links to a.md, same as column A.
$wsZh.Hyperlinks.Add($wsZh.Range("J2"), $aMdUrl, [Type]::Missing, [Type]::Missing, "a.md")
$wsZh.Hyperlinks.Add($wsZh.Range("J3"), $aMdUrl, [Type]::Missing, [Type]::Missing, "a.md")

# "Latest Handback File" (K) and "Latest Handback DateTime" (L)
$wsZh.Range("K2").Value = $zhHandback
$wsZh.Range("K3").Value = $zhHandback
$wsZh.Range("L2").Value = $zhDateTime
$wsZh.Range("L3").Value = $zhDateTime

$wsZh.Columns.Item(3).ColumnWidth = 29.14
$wsZh.Columns.Item(11).ColumnWidth = 39.17

# ---------------------------------------------------------------------
# de-de sheet: record the handback for a.md
# ---------------------------------------------------------------------
$wsDe = $wb.Worksheets.Item("de-de")

$wsDe.Range("C2").Value = $newStatus
$wsDe.Range("C3").Value = $newStatus

$wsDe.Hyperlinks.Add($wsDe.Range("J2"), $aMdUrl, [Type]::Missing, [Type]::Missing, "a.md")
$wsDe.Hyperlinks.Add($wsDe.Range("J3"), $aMdUrl, [Type]::Missing, [Type]::Missing, "a.md")

$wsDe.Range("K2").Value = $deHandback
$wsDe.Range("K3").Value = $deHandback
$wsDe.Range("L2").Value = $deDateTime
$wsDe.Range("L3").Value = $deDateTime

$wsDe.Columns.Item(3).ColumnWidth = 29.14
$wsDe.Columns.Item(11).ColumnWidth = 39.17
